# Word COM-interop script applying the tracked changes:
#   1. Split the "left-bottom" run into "left-" + "middle" (the new
#      "middle" run drops the explicit black color, keeping only rtl=0).
#   2. Remove the comment anchored around the "p" of "panier" (both the
#      range markers/reference in the body and the comment definition).
#   3. Give the section an explicit footer distance (w:footer="720").

$d = $word.ActiveDocument

# --- 1. "left-bottom" -> "left-" / "middle" -------------------------------
$rng = $d.Content
$rng.Find.Execute("left-bottom", $true, $false, $false, $false, $false,
                   $true, 1, $false, $null, 0)
$matchStart = $rng.Start
$matchEnd = $rng.End

# Give the "bottom" tail automatic color *before* retyping it so it no
# longer matches the preceding "left-" run's direct black color -- that
# keeps it from being merged back into a single run once the text changes.
$tail = $d.Range($matchStart + 5, $matchEnd)
$tail.Font.Color = -16777216
$tail.Text = "middle"

# --- 2. Remove the comment on "p" of "panier" -----------------------------
$d.Comments(1).Delete()

# --- 3. Explicit footer distance on the section's page margins -----------
$d.Sections(1).PageSetup.FooterDistance = 36
